$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -12.672
    6   = -12.239
    7   = -12.89
    8   = -12.822
    16  = -13.093
    20  = -12.13
    21  = -12.725
    28  = -12.766
    29  = -12.129
    30  = -12.202
    32  = -12.628
    40  = -12.446
    46  = -13.654
    51  = -11.296
    52  = -11.382
    57  = -13.818
    59  = -12.803
    62  = -13.211
    66  = -11.713
    73  = -12.298
    74  = -11.982
    77  = -12.596
    92  = -10.753
    100 = -12.221
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
